$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "60.862.66"
$ws.Range("E2").Value = "  -3.21%  "

$ws.Range("D3").Value = "2.919.41"
$ws.Range("E3").Value = "  -3.82%  "

$ws.Range("E4").Value = "  +0.09%  "

$ws.Range("D5").Value = "588.13"
$ws.Range("E5").Value = "  -1.26%  "

$ws.Range("D6").Value = "146.48"
$ws.Range("E6").Value = "  -4.63%  "

$ws.Range("E7").Value = "  +0.09%  "

$ws.Range("B8").Value = "LidoStakedEther"
$ws.Range("C8").Value = "https://coinranking.com/coin/VINVMYf0u+lidostakedether-steth"
$ws.Range("D8").Value = "2.919.79"
$ws.Range("E8").Value = "  -3.73%  "

$ws.Range("B9").Value = "XRP"
$ws.Range("C9").Value = "https://coinranking.com/coin/-l8Mn2pVlRs-p+xrp-xrp"
$ws.Range("D9").Value = "0.503"
$ws.Range("E9").Value = "  -2.94%  "

$ws.Range("D10").Value = "6.79"
$ws.Range("E10").Value = "  +6.67%  "

$ws.Range("D11").Value = "0.145"
$ws.Range("E11").Value = "  -4.20%  "

$ws.Range("D12").Value = "0.447"
$ws.Range("E12").Value = "  -3.76%  "

$ws.Range("D13").Value = "0.0000226"
$ws.Range("E13").Value = "  -3.42%  "

$ws.Range("D14").Value = "33.71"
$ws.Range("E14").Value = "  -4.45%  "

$ws.Range("D15").Value = "0.127"
$ws.Range("E15").Value = "  +0.64%  "

$ws.Range("D16").Value = "3.404.74"
$ws.Range("E16").Value = "  -3.73%  "

$ws.Range("D17").Value = "60.849.51"
$ws.Range("E17").Value = "  -3.17%  "

$ws.Range("D18").Value = "6.75"
$ws.Range("E18").Value = "  -4.68%  "

$ws.Range("D19").Value = "2.926.07"
$ws.Range("E19").Value = "  -3.47%  "

$ws.Range("D20").Value = "429.70"
$ws.Range("E20").Value = "  -5.31%  "

$ws.Range("D21").Value = "13.68"
$ws.Range("E21").Value = "  -4.32%  "

$ws.Range("D22").Value = "0.681"
$ws.Range("E22").Value = "  -2.04%  "

$ws.Range("D23").Value = "7.16"
$ws.Range("E23").Value = "  -4.87%  "

$ws.Range("B24").Value = "RenderToken"
$ws.Range("C24").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D24").Value = "11.13"
$ws.Range("E24").Value = "  +1.22%  "

$ws.Range("B25").Value = "Litecoin"
$ws.Range("C25").Value = "https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc"
$ws.Range("D25").Value = "80.67"
$ws.Range("E25").Value = "  -3.21%  "

$ws.Range("D26").Value = "2.25"
$ws.Range("E26").Value = "  -2.08%  "

$ws.Range("D27").Value = "11.90"
$ws.Range("E27").Value = "  -3.20%  "

$ws.Range("E28").Value = "  -0.09%  "

$ws.Range("E29").Value = "  +0.15%  "

$ws.Range("D30").Value = "7.26"
$ws.Range("E30").Value = "  -3.59%  "

$ws.Range("D31").Value = "2.63"
$ws.Range("E31").Value = "  -2.68%  "

$ws.Range("D32").Value = "2.18"
$ws.Range("E32").Value = "  -2.29%  "

$ws.Range("D33").Value = "26.63"
$ws.Range("E33").Value = "  -4.01%  "

$ws.Range("D34").Value = "0.107"
$ws.Range("E34").Value = "  -3.41%  "

$ws.Range("D35").Value = "0.0₃0880"
$ws.Range("E35").Value = "  +3.01%  "

$ws.Range("E36").Value = "  -2.80%  "

$ws.Range("D37").Value = "5.66"
$ws.Range("E37").Value = "  -4.59%  "

$ws.Range("D38").Value = "3.04"
$ws.Range("E38").Value = "  -2.95%  "

$ws.Range("B39").Value = "Kaspa"
$ws.Range("C39").Value = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
$ws.Range("D39").Value = "0.128"
$ws.Range("E39").Value = "  +1.62%  "

$ws.Range("B40").Value = "OKB"
$ws.Range("C40").Value = "https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb"
$ws.Range("D40").Value = "49.82"
$ws.Range("E40").Value = "  -1.79%  "

$ws.Range("D41").Value = "2.01"
$ws.Range("E41").Value = "  -3.56%  "

$ws.Range("D42").Value = "8.71"
$ws.Range("E42").Value = "  -4.69%  "

$ws.Range("D43").Value = "0.301"
$ws.Range("E43").Value = "  -1.00%  "

$ws.Range("D44").Value = "42.07"
$ws.Range("E44").Value = "  +0.10%  "

$ws.Range("D45").Value = "380.15"
$ws.Range("E45").Value = "  -3.62%  "

$ws.Range("D46").Value = "0.0350"
$ws.Range("E46").Value = "  -2.65%  "

$ws.Range("D47").Value = "2.675.73"
$ws.Range("E47").Value = "  -2.56%  "

$ws.Range("D48").Value = "131.89"
$ws.Range("E48").Value = "  -0.18%  "

$ws.Range("E49").Value = "  +0.01%  "

$ws.Range("D50").Value = "24.79"
$ws.Range("E50").Value = "  +2.64%  "

$ws.Range("E51").Value = "  -2.06%  "
